# Update loading-times data after modifying penalties/fixing assignments
# and generating a new population. This rewrites the data rows (2-11) and
# appends a new row (12) to reflect the new simulation results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5, 1, 5, 5),
    @(10, 1, 10, 12),
    @(3, 3, 5, 5),
    @(2, 4, 5, 5),
    @(1, 5, 5, 5),
    @(6, 5, 10, 10),
    @(8, 5, 15, 15),
    @(9, 5, 20, 20),
    @(10, 5, 25, 27),
    @(4, 6, 5, 5),
    @(7, 6, 10, 10)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
